$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the abstract text in C3 (typo fix: removed stray space before "was done.") ---
$ws.Range("C3").Value = 'abstract for the main project "AUTOMATIC RESUME CORRECTION"was done.'

# --- New daily-routine rows (5-8) ---

# Row 5: placement preparation
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 43597
$ws.Range("B5").Value = "Prepared for placement(quantitative,reasoning and verbal aptitude)"

# Row 6: modules / prototype
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 43720
$ws.Range("B6").Value = "modules were added,prototype for the project was drawn"
$ws.Range("C6").Value = "modules were designed."

# Row 7: web pages / register-login design
$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 43750
$ws.Range("B7").Value = "register,login page were designed. The overall website and flow of it was drawn."
$ws.Range("C7").Value = "web pages drawn related to the project."

# Row 8: login and home page design
$ws.Range("A4").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 43781
$ws.Range("B8").Value = "login and home page were designed."

# --- Match the final selection state from the workbook ---
$ws.Range("C8").Select()
